# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback/localization round has completed (the two xliff files have been
# returned and are now in sync with en-US), filling in the "Latest Target
# File", "Latest Handback File" and "Latest Handback DateTime" columns on
# the zh-cn and de-de sheets, adding hyperlinks for the newly filled in
# "Latest Target File" cells, updating the status text on the Overview
# sheet, and widening a few columns so the longer file names are readable.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item(1)
$ZhCn     = $wb.Worksheets.Item(2)
$DeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Overview sheet: status changes from "In Translation" to "Handed
#    back: in sync with en-US" for both locale columns / both rows.
# ---------------------------------------------------------------------
$Overview.Range("E2").Value = "Handed back: in sync with en-US"
$Overview.Range("F2").Value = "Handed back: in sync with en-US"
$Overview.Range("E3").Value = "Handed back: in sync with en-US"
$Overview.Range("F3").Value = "Handed back: in sync with en-US"

# Existing hyperlink cells (column A) use an underlined, cornflower-blue
# (#6495ED) font. Re-create that exact look on the newly-populated
# "Latest Target File" cells (I2/I3) by setting the font directly,
# rather than via a named style, so it matches the existing look
# instead of introducing a second, differently-coloured "hyperlink"
# style. xlUnderlineStyleSingle = 2; the Color is a BGR-packed long for
# RGB(0x64,0x95,0xED).
$HyperlinkUnderline = 2
$HyperlinkColor = 15570276

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I), Latest Handback File
#    (J) and Latest Handback DateTime (K) for both data rows, and turn
#    the newly-populated target-file cells into hyperlinks (matching
#    the existing hyperlink on column A for the same row).
# ---------------------------------------------------------------------
$ZhCn.Range("I2").Value = "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md"
$ZhCn.Range("I2").Font.Underline = $HyperlinkUnderline
$ZhCn.Range("I2").Font.Color = $HyperlinkColor
$ZhCn.Range("J2").Value = "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.fff98a315c0ee0dc2e3812f8c60d3a9663c66ad5.zh-cn.xlf"
$ZhCn.Range("K2").Value = "2016-08-25 10:25:51"

$ZhCn.Range("I3").Value = "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md"
$ZhCn.Range("I3").Font.Underline = $HyperlinkUnderline
$ZhCn.Range("I3").Font.Color = $HyperlinkColor
$ZhCn.Range("J3").Value = "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.c9fbe2d75ae13341fdb09529e0b462a361ca90b6.zh-cn.xlf"
$ZhCn.Range("K3").Value = "2016-08-25 10:25:51"

$ZhCn.Hyperlinks.Add($ZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a41474a6f5b4429217ef36d038b270789e520c3/e2e/58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md", "", "", "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md")
$ZhCn.Hyperlinks.Add($ZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a41474a6f5b4429217ef36d038b270789e520c3/e2e/d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md", "", "", "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md")

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback datetime.
# ---------------------------------------------------------------------
$DeDe.Range("I2").Value = "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md"
$DeDe.Range("I2").Font.Underline = $HyperlinkUnderline
$DeDe.Range("I2").Font.Color = $HyperlinkColor
$DeDe.Range("J2").Value = "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.fff98a315c0ee0dc2e3812f8c60d3a9663c66ad5.de-de.xlf"
$DeDe.Range("K2").Value = "2016-08-25 10:25:59"

$DeDe.Range("I3").Value = "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md"
$DeDe.Range("I3").Font.Underline = $HyperlinkUnderline
$DeDe.Range("I3").Font.Color = $HyperlinkColor
$DeDe.Range("J3").Value = "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.c9fbe2d75ae13341fdb09529e0b462a361ca90b6.de-de.xlf"
$DeDe.Range("K3").Value = "2016-08-25 10:25:59"

$DeDe.Hyperlinks.Add($DeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a41474a6f5b4429217ef36d038b270789e520c3/e2e/58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md", "", "", "58dbcf0a-3ba8-4948-ac2b-4a3b741bbe61.md")
$DeDe.Hyperlinks.Add($DeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a41474a6f5b4429217ef36d038b270789e520c3/e2e/d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md", "", "", "d6a5cd32-a46a-4c5d-b69e-d49cc58f9a9d.md")

# ---------------------------------------------------------------------
# 4. Widen columns so the newly-populated / longer values are readable.
#    (ColumnWidth is expressed in characters and gets pixel-snapped by
#    the engine just like real Excel; the inputs below are chosen so the
#    stored width lands on the intended value.)
# ---------------------------------------------------------------------
$Overview.Columns.Item(5).ColumnWidth = 29.166666666666664   # E -> ~29.98
$Overview.Columns.Item(6).ColumnWidth = 29.166666666666664   # F -> ~29.98

$ZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666664      # C -> ~29.98
$ZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I -> 40
$ZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J -> 40

$DeDe.Columns.Item(3).ColumnWidth  = 29.166666666666664      # C -> ~29.98
$DeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I -> 40
$DeDe.Columns.Item(10).ColumnWidth = 39.166666666666664      # J -> 40

Write-Host "Handback report generated."
